# Applies the "add lulucf and uba data 2023" edit:
#  - Fills in the 2024-11 (column DH) figures for rows 2-46, which were
#    previously blank.
#  - For the two "Fossil energy [FE]" rows (45 = GWH, 46 = PC) the monthly
#    history in columns R:CW (2017-01 .. 2023-12) is cleared out as part of
#    the same data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New DH (2024-11) values, keyed by row number.
$dhValues = @{
    2  = 5912.698
    3  = 2037.953
    4  = 34.467
    5  = 228.079
    6  = 3.857
    7  = 60.176
    8  = 1.018
    9  = 116.7
    10 = 1.974
    11 = 1583.811
    12 = 26.787
    13 = 49.187
    14 = 0.832
    15 = 3796.951
    16 = 64.217
    17 = 2849.504
    18 = 48.193
    19 = 1627.195
    20 = 27.52
    21 = 1222.309
    22 = 20.673
    23 = 0
    24 = 0
    25 = 0
    26 = 0
    27 = 719.3680000000001
    28 = 12.166
    29 = 719.3680000000001
    30 = 12.166
    31 = 0
    32 = 0
    33 = 0
    34 = 0
    35 = 0
    36 = 0
    37 = 0
    38 = 0
    39 = 0
    40 = 0
    41 = 0
    42 = 0
    43 = 305.873
    44 = 5.173
    45 = 1809.874
    46 = 30.61
}

foreach ($row in $dhValues.Keys) {
    $ws.Range("DH$row").Value = $dhValues[$row]
}

# Rows 45 and 46 ("Fossil energy [FE]", GWH and PC) had their 2017-01
# through 2023-12 monthly values (columns R:CW) cleared.
$ws.Range("R45:CW45").Value = ""
$ws.Range("R46:CW46").Value = ""
